$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 258, pushing the existing
# rows 258:281 down to 260:283.
$ws.Rows.Item(258).Insert()
$ws.Rows.Item(258).Insert()

# New row 258 data
$ws.Cells.Item(258, 1).Value = 6
$ws.Cells.Item(258, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(258, 3).Value = "Metropolitana"
$ws.Cells.Item(258, 4).Value = 44783
$ws.Cells.Item(258, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(258, 5).Value = 13
$ws.Cells.Item(258, 6).Value = 100112026
$ws.Cells.Item(258, 7).Value = "Haba"
$ws.Cells.Item(258, 8).Value = "Sin especificar"
$ws.Cells.Item(258, 9).Value = "Primera"
$ws.Cells.Item(258, 10).Value = 230
$ws.Cells.Item(258, 11).Value = 12000
$ws.Cells.Item(258, 12).Value = 12000
$ws.Cells.Item(258, 13).Value = 12000
$ws.Cells.Item(258, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(258, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(258, 16).Value = 480
$ws.Cells.Item(258, 17).Value = 25
$ws.Cells.Item(258, 18).Value = "Hortaliza"

# New row 259 data
$ws.Cells.Item(259, 1).Value = 6
$ws.Cells.Item(259, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(259, 3).Value = "Metropolitana"
$ws.Cells.Item(259, 4).Value = 44783
$ws.Cells.Item(259, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(259, 5).Value = 13
$ws.Cells.Item(259, 6).Value = 100112026
$ws.Cells.Item(259, 7).Value = "Haba"
$ws.Cells.Item(259, 8).Value = "Sin especificar"
$ws.Cells.Item(259, 9).Value = "Primera"
$ws.Cells.Item(259, 10).Value = 170
$ws.Cells.Item(259, 11).Value = 14000
$ws.Cells.Item(259, 12).Value = 14000
$ws.Cells.Item(259, 13).Value = 14000
$ws.Cells.Item(259, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(259, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(259, 16).Value = 560
$ws.Cells.Item(259, 17).Value = 25
$ws.Cells.Item(259, 18).Value = "Hortaliza"
